$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 ("Полученная частота
# дискретизации аудио" and below), shifting the CIC-compensator rows
# down by one.
$ws.Rows.Item(10).Insert()

# New row 10: sampling rate obtained straight from the CIC compensator
# (before dividing by the compensator decimation factor C8/D8).
$ws.Range("B10").Value = "Полученная частота дискретизации CIC коспенсатора"
$ws.Range("C10").Formula = "=C7/C9"
$ws.Range("C10").NumberFormat = "0.00"
$ws.Range("D10").Formula = "=D7/D9"

# Center the "Практика" value header over column D.
$ws.Range("D6").Font.Bold = $true
$ws.Range("D6").HorizontalAlignment = -4108

# Updated clock-rate-at-CIC-compensator input value.
$ws.Range("C4").Value = 34

# Restore the user selection to the edited input cell.
$ws.Range("C4").Select() | Out-Null
